$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = "tfa@SS01"
$ws.Range("F10").Value = "tfa@SS01"
$ws.Range("B11").Value = "Est. Engenharia"

$ws.Range("A13").Value = "TFA"
$ws.Range("B13").Value = "Máq. Da manutenção"
$ws.Range("C13").Value = "10.28.2.124"
$ws.Range("D13").Value = "Hostname"
$ws.Range("E13").Value = "TATFASS03\VALE"
$ws.Range("F13").Value = "V@le#tatf@02"

$ws.Range("A14").Value = "TFA"
$ws.Range("B14").Value = "Máq. Da manutenção"
$ws.Range("C14").Value = "10.28.2.124"
$ws.Range("D14").Value = "Hostname"
$ws.Range("E14").Value = "TATFASS03\Manut-TFA"
$ws.Range("F14").Value = "M@nutencao"
